{"js": "// Update the date paragraph (first paragraph of the body) and the\n// two-digit-number / one-digit-number division problems laid out in the\n// first table of the document. Cells are addressed positionally (row /\n// column index) rather than by matching old text, because several problem\n// cells share the same original text (e.g. \"29\u00f78=\" appears twice) but map\n// to different new values.\n\n// 1) Update the date line: \"2025-06-21 Saturday\" -> \"2025-06-22 Sunday\"\nconst dateResults = context.document.body.search(\"2025-06-21 Saturday\", { matchCase: true });\ndateResults.load(\"items\");\nawait context.sync();\n\nfor (const r of dateResults.items) {\n  r.insertText(\"2025-06-22 Sunday\", \"Replace\");\n}\nawait context.sync();\n\n// 2) Update the division problems inside the first table, addressed by\n//    (rowIndex, cellIndex) so duplicate original values don't collide.\nconst table = context.document.body.tables.getFirst();\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// rowIndex -> [newValue for cell0, cell1, cell2, cell3, cell4]\nconst rowUpdates = {\n  0: [\"42\u00f74=\", \"94\u00f79=\", \"15\u00f76=\", \"98\u00f74=\", \"68\u00f79=\"],\n  4: [\"40\u00f75=\", \"62\u00f72=\", \"60\u00f73=\", \"12\u00f79=\", \"85\u00f73=\"],\n  8: [\"37\u00f76=\", \"26\u00f78=\", \"51\u00f72=\", \"33\u00f72=\", \"71\u00f72=\"],\n  12: [\"80\u00f76=\", \"91\u00f76=\", \"50\u00f76=\", \"71\u00f72=\", \"32\u00f78=\"],\n  16: [\"79\u00f72=\", \"64\u00f74=\", \"98\u00f78=\", \"28\u00f74=\", \"23\u00f76=\"],\n};\n\nfor (const rowIndexStr of Object.keys(rowUpdates)) {\n  const rowIndex = Number(rowIndexStr);\n  const newValues = rowUpdates[rowIndex];\n  const cells = rows.items[rowIndex].cells;\n  cells.load(\"items\");\n  await context.sync();\n  for (let c = 0; c < newValues.length; c++) {\n    cells.items[c].value = newValues[c];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date paragraph and the two-digit-number / one-digit-number\n# division problems laid out in the first table of the document.\n#\n# Table cells are addressed positionally (row / column index) rather than\n# by matching old text, because several problem cells share the same\n# original text (e.g. \"29\u00f78=\" appears twice) but map to different new\n# values.\n\n$d = $word.ActiveDocument\n\n# 1) Update the date line: \"2025-06-21 Saturday\" -> \"2025-06-22 Sunday\"\n$findRange = $d.Content\n$findRange.Find.Execute(\"2025-06-21 Saturday\", $false, $false, $false, $false, $false, $true, 1, $false, \"2025-06-22 Sunday\", 2) | Out-Null\n\n# 2) Update the division problems inside the first table, addressed by\n#    (row, column) so duplicate original values don't collide.\n$table = $d.Tables.Item(1)\n\n# 1-based row index -> new values for columns 1..5\n$rowUpdates = @{\n    1  = @(\"42\u00f74=\", \"94\u00f79=\", \"15\u00f76=\", \"98\u00f74=\", \"68\u00f79=\")\n    5  = @(\"40\u00f75=\", \"62\u00f72=\", \"60\u00f73=\", \"12\u00f79=\", \"85\u00f73=\")\n    9  = @(\"37\u00f76=\", \"26\u00f78=\", \"51\u00f72=\", \"33\u00f72=\", \"71\u00f72=\")\n    13 = @(\"80\u00f76=\", \"91\u00f76=\", \"50\u00f76=\", \"71\u00f72=\", \"32\u00f78=\")\n    17 = @(\"79\u00f72=\", \"64\u00f74=\", \"98\u00f78=\", \"28\u00f74=\", \"23\u00f76=\")\n}\n\nforeach ($rowIndex in $rowUpdates.Keys) {\n    $newValues = $rowUpdates[$rowIndex]\n    for ($col = 1; $col -le $newValues.Length; $col++) {\n        $cell = $table.Cell($rowIndex, $col)\n        $cell.Range.Text = $newValues[$col - 1]\n    }\n}\n"}
